{"js": "// Insert a new bulleted list item \"Not uploading? Unplug power supply\"\n// right after the \"Maybe try run bootloader.\" item (end of the\n// \"Programming\" list), matching the target XML diff.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"Maybe try run bootloader.\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(`Could not find anchor paragraph: \"${anchorText}\"`);\n}\n\n// Inserting a new paragraph directly after an existing list paragraph\n// inherits that paragraph's style/list numbering (same as pressing Enter\n// at the end of the line in Word), producing the same ListParagraph /\n// numId=2 bullet formatting used by the diff.\nanchor.insertParagraph(\"Not uploading? Unplug power supply\", \"After\");\n\nawait context.sync();\n", "ps1": "# Insert a new bulleted list item \"Not uploading? Unplug power supply\"\n# right after the \"Maybe try run bootloader.\" item (end of the\n# \"Programming\" list), matching the target XML diff.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"Maybe try run bootloader.\"\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq $anchorText) {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not find anchor paragraph: $anchorText\"\n}\n\n# Appending \"<CR>text\" to the end of the anchor paragraph's range creates a\n# new paragraph right after it that inherits the anchor's style/list\n# numbering (same as pressing Enter at the end of the line in Word),\n# producing the same ListParagraph / numId=2 bullet formatting used by the\n# diff.\n$target.Range.InsertAfter([char]13 + \"Not uploading? Unplug power supply\")\n"}
